# The authored commit swaps the contents of ppt/theme/theme1.xml (the
# slide-master theme, originally the "Integral" colour scheme) and
# ppt/theme/theme2.xml (the notes-master theme, originally the stock
# "Office Theme" colour scheme) - i.e. after the edit the slide master
# uses the "Office" palette and the notes master would use "Integral".
#
# The PowerPoint object model only exposes the slide-master's theme
# (ActivePresentation.SlideMaster.Theme); the notes-master theme is not
# independently reachable through COM automation (NotesMaster resolves
# back to the SlideMaster object), so we apply the reachable half of
# the swap: re-point the slide master's 12-slot theme colour scheme to
# the "Office Theme" palette that used to live in theme2.xml.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# RGB(r,g,b) COM colour integer = r + g*256 + b*65536
# Order is: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$tcs.Item(1).RGB  = 0x000000 + 0x00*256 + 0x00*65536   # dk1    000000
$tcs.Item(2).RGB  = 0xFF + 0xFF*256 + 0xFF*65536        # lt1    FFFFFF
$tcs.Item(3).RGB  = 0x44 + 0x54*256 + 0x6A*65536        # dk2    44546A
$tcs.Item(4).RGB  = 0xE7 + 0xE6*256 + 0xE6*65536        # lt2    E7E6E6
$tcs.Item(5).RGB  = 0x5B + 0x9B*256 + 0xD5*65536        # accent1 5B9BD5
$tcs.Item(6).RGB  = 0xED + 0x7D*256 + 0x31*65536        # accent2 ED7D31
$tcs.Item(7).RGB  = 0xA5 + 0xA5*256 + 0xA5*65536        # accent3 A5A5A5
$tcs.Item(8).RGB  = 0xFF + 0xC0*256 + 0x00*65536        # accent4 FFC000
$tcs.Item(9).RGB  = 0x44 + 0x72*256 + 0xC4*65536        # accent5 4472C4
$tcs.Item(10).RGB = 0x70 + 0xAD*256 + 0x47*65536        # accent6 70AD47
$tcs.Item(11).RGB = 0x05 + 0x63*256 + 0xC1*65536        # hlink   0563C1
$tcs.Item(12).RGB = 0x95 + 0x4F*256 + 0x72*65536        # folHlink 954F72
